# Apply the Xhosa translation edits described in the commit diff.
$d = $word.ActiveDocument

# 1) Title paragraph: "Translated" -> "Iphepha loku Bhalisa " (keep "Text_V3" run,
#    commentRangeStart/End and commentReference untouched). The word "Translated"
#    is exactly the first 10 characters of the document.
$d.Range(0, 10).Text = "Iphepha loku Bhalisa "

# 2) "Ingaba ungumzali okanye umnonopheli womntwana?" -> "Ingaba ungumgcini womntwana?"
$d.Paragraphs.Item(4).Range.Text = "Ingaba ungumgcini womntwana?"

# 3) First "translated" placeholder -> "Phuhlisa ubuzali bakho ngeentsuku ezi 10!"
$d.Paragraphs.Item(6).Range.Text = "Phuhlisa ubuzali bakho ngeentsuku ezi 10!"

# 4) Second "translated" placeholder -> WhatsApp instructions line
$d.Paragraphs.Item(8).Range.Text = "WattsAppa {'Molo'} ku {0....} kwaye ulandele umgaqo ukuzibandakanya nenkqubo ye {ParentText}"

# 5) Third "translated" placeholder (had a trailing line break in the original run)
#    -> website/contact details line, with no trailing break.
$d.Paragraphs.Item(10).Range.Text = "Ngenkcukacha ezithe vetshe kunye nemibuzo, jonga i webhisayithi yethu kwi {www.globalparenting.org} okanye imeyile {swift@globalparenting.org}"

# 6) Update the reviewer comment text.
$d.Comments.Item(1).Content = "@chiara.facciola@idems.international Nantsi i imbalelwano ye phepha le posta. Onke amagama abiyelwe nge brakhethi ahlale engaguqulwanga. Sizokufaka inombolo ka wattsapp xana sithe sanayo. Enkosi"
